$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

# First, copy the bold/border style used by column A (rows 2-6) down to the
# two new rows (7-8) that this edit introduces, so they match the existing
# label formatting.
$ws.Range("A6").Copy() | Out-Null
$ws.Range("A7:A8").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 2: start -> page_width
$ws.Range("A2").Value = "page_width"
$ws.Range("B2").Value = 200
$ws.Range("C2").Value = "The width of the page in arbitrary units"

# Row 3: resolution -> page_height
$ws.Range("A3").Value = "page_height"
$ws.Range("B3").Value = 100
$ws.Range("C3").Value = "The height of the page in arbitrary units"

# Row 4: width -> start_date
$ws.Range("A4").Value = "start_date"
# Force the date-like string to be stored as literal text (not an Excel
# date serial number), matching the original "2024-01-01" inline string.
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "2024-01-01"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = "The start date of the chart"

# Row 5: height -> header_height
$ws.Range("A5").Value = "header_height"
$ws.Range("B5").Value = 0.1
$ws.Range("C5").Value = "The height of the header as a proportion of the page height"

# Row 6: rows -> footer_height
$ws.Range("A6").Value = "footer_height"
$ws.Range("B6").Value = 0.1
$ws.Range("C6").Value = "The height of the footer as a proportion of the page height"

# Row 7 (new): total_scale_height
$ws.Range("A7").Value = "total_scale_height"
$ws.Range("B7").Value = 0.1
$ws.Range("C7").Value = "The total height available for scales as a proportion of the page height"

# Row 8 (new): row_quantity
$ws.Range("A8").Value = "row_quantity"
$ws.Range("B8").Value = 5
$ws.Range("C8").Value = "The number of rows in the chart"
